$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the input values that drive the formulas on the sheet.
$ws.Range("A2").Value = 1.95
$ws.Range("C2").Value = 27000
$ws.Range("D2").Value = 30000

$ws.Range("A5").Value = 680
$ws.Range("A6").Value = 680
$ws.Range("A7").Value = 680

# Recalculate so dependent formula cells pick up the new cached values.
$excel.Calculate()

# Reflect the author's final cell selection on the sheet.
$ws.Range("A8").Select()
